# ---------------------------------------------------------------------------
# SeniorConnect master log - append freshly-captured sensor readings
#
# On 2026-01-30 (~16:57-16:59) several new sensor events were logged.
# This script appends the corresponding rows to the bottom of the six
# affected sheets (ALERTS, PIR, Humidity, Temperature, Proximity, mmWave).
# Every sheet shares the same column layout:
#   A=Date  B=Timestamp  C=Hour  D=Location  E=Value  F=Status
# Only brand-new rows (after each sheet's current last row) are written;
# no existing cell is modified.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Writes one data row (Date, Timestamp, Hour, Location, Value, Status) into
# row $RowNum of worksheet $Sheet. Values that look like dates (column A) or
# percentages (column E) are forced to Text format first, because Excel's
# COM layer otherwise silently re-interprets strings such as "2026-01-30"
# or "86.8%" as a date serial number / numeric percentage rather than
# keeping them as the plain text used throughout the rest of the log.
function Add-LogRow {
    param($Sheet, $RowNum, $DateStr, $TimeStr, $HourStr, $Location, $Value, $Status)

    $Sheet.Cells.Item($RowNum, 1).NumberFormat = "@"
    $Sheet.Cells.Item($RowNum, 1).Value = $DateStr
    $Sheet.Cells.Item($RowNum, 2).Value = $TimeStr
    $Sheet.Cells.Item($RowNum, 3).Value = $HourStr
    $Sheet.Cells.Item($RowNum, 4).Value = $Location

    if ($Value -like "*%*") {
        $Sheet.Cells.Item($RowNum, 5).NumberFormat = "@"
    }
    $Sheet.Cells.Item($RowNum, 5).Value = $Value
    $Sheet.Cells.Item($RowNum, 6).Value = $Status
}

# --- "ALERTS" (sheet 1) -- append rows 9-10 ---
$ws = $wb.Worksheets.Item(1)
Add-LogRow $ws 9 "2026-01-30" "16:58:04" "16:00" "Living Room" "CRITICAL EMERGENCY" "FALL_DETECTED"
Add-LogRow $ws 10 "2026-01-30" "16:58:09" "16:00" "Living Room" "CRITICAL EMERGENCY" "FALL_DETECTED"

# --- "PIR" (sheet 2) -- append rows 101-116 ---
$ws = $wb.Worksheets.Item(2)
Add-LogRow $ws 101 "2026-01-30" "16:57:14" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 102 "2026-01-30" "16:57:15" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 103 "2026-01-30" "16:57:18" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 104 "2026-01-30" "16:58:09" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 105 "2026-01-30" "16:58:10" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 106 "2026-01-30" "16:58:14" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 107 "2026-01-30" "16:58:19" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 108 "2026-01-30" "16:58:23" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 109 "2026-01-30" "16:58:28" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 110 "2026-01-30" "16:58:34" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 111 "2026-01-30" "16:58:38" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 112 "2026-01-30" "16:58:44" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 113 "2026-01-30" "16:58:44" "16:00" "Living Room" "RECOVERY_DETECTION" "Inactive"
Add-LogRow $ws 114 "2026-01-30" "16:58:48" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 115 "2026-01-30" "16:58:54" "16:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 116 "2026-01-30" "16:58:59" "16:00" "Bathroom" "No Motion" "Inactive"

# --- "Humidity" (sheet 3) -- append rows 63-73 ---
$ws = $wb.Worksheets.Item(3)
Add-LogRow $ws 63 "2026-01-30" "16:57:15" "16:00" "Bathroom" "86.8%" "Active"
Add-LogRow $ws 64 "2026-01-30" "16:57:16" "16:00" "Bathroom" "87.6%" "Active"
Add-LogRow $ws 65 "2026-01-30" "16:57:18" "16:00" "Bathroom" "86.7%" "Active"
Add-LogRow $ws 66 "2026-01-30" "16:58:09" "16:00" "Bathroom" "87.6%" "Active"
Add-LogRow $ws 67 "2026-01-30" "16:58:19" "16:00" "Bathroom" "87.6%" "Active"
Add-LogRow $ws 68 "2026-01-30" "16:58:24" "16:00" "Bathroom" "87.6%" "Active"
Add-LogRow $ws 69 "2026-01-30" "16:58:29" "16:00" "Bathroom" "87.6%" "Active"
Add-LogRow $ws 70 "2026-01-30" "16:58:34" "16:00" "Bathroom" "87.6%" "Active"
Add-LogRow $ws 71 "2026-01-30" "16:58:49" "16:00" "Bathroom" "87.6%" "Active"
Add-LogRow $ws 72 "2026-01-30" "16:58:55" "16:00" "Bathroom" "87.6%" "Active"
Add-LogRow $ws 73 "2026-01-30" "16:58:59" "16:00" "Bathroom" "87.6%" "Active"

# --- "Temperature" (sheet 4) -- append rows 25-35 ---
$ws = $wb.Worksheets.Item(4)
Add-LogRow $ws 25 "2026-01-30" "16:57:15" "16:00" "Bathroom" "22.7C" "Active"
Add-LogRow $ws 26 "2026-01-30" "16:57:16" "16:00" "Bathroom" "22.6C" "Active"
Add-LogRow $ws 27 "2026-01-30" "16:57:19" "16:00" "Bathroom" "22.6C" "Active"
Add-LogRow $ws 28 "2026-01-30" "16:58:10" "16:00" "Bathroom" "22.6C" "Active"
Add-LogRow $ws 29 "2026-01-30" "16:58:19" "16:00" "Bathroom" "22.6C" "Active"
Add-LogRow $ws 30 "2026-01-30" "16:58:24" "16:00" "Bathroom" "22.7C" "Active"
Add-LogRow $ws 31 "2026-01-30" "16:58:29" "16:00" "Bathroom" "22.6C" "Active"
Add-LogRow $ws 32 "2026-01-30" "16:58:34" "16:00" "Bathroom" "22.7C" "Active"
Add-LogRow $ws 33 "2026-01-30" "16:58:49" "16:00" "Bathroom" "22.7C" "Active"
Add-LogRow $ws 34 "2026-01-30" "16:58:55" "16:00" "Bathroom" "22.7C" "Active"
Add-LogRow $ws 35 "2026-01-30" "16:58:59" "16:00" "Bathroom" "22.6C" "Active"

# --- "Proximity" (sheet 5) -- append rows 32-33 ---
$ws = $wb.Worksheets.Item(5)
Add-LogRow $ws 32 "2026-01-30" "16:59:00" "16:00" "Living Room Main Door" "Clear" "Inactive"
Add-LogRow $ws 33 "2026-01-30" "16:59:01" "16:00" "Living Room Main Door" "Detected" "Active"

# --- "mmWave" (sheet 6) -- append rows 25-26 ---
$ws = $wb.Worksheets.Item(6)
Add-LogRow $ws 25 "2026-01-30" "16:58:44" "16:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $ws 26 "2026-01-30" "16:58:55" "16:00" "Living Room" "PRESENCE_DETECTED" "Active"

